$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: add ml_analyzed/runtime/inhibittime/numtriggers values, fix Volume to "2*" ---
$ws.Range("B2").Value = "2*"
$ws.Range("E2").Value = 1.6185
$ws.Range("F2").Value = 478
$ws.Range("G2").Value = 73.459999999999994
$ws.Range("H2").Value = 1005

# --- Row 3: add ml_analyzed/runtime/inhibittime/numtriggers values ---
$ws.Range("E3").Value = 1.637
$ws.Range("F3").Value = 478
$ws.Range("G3").Value = 69.22
$ws.Range("H3").Value = 957

# --- Row 4: fix Volume to "5*", add ml_analyzed/runtime/inhibittime/numtriggers values ---
$ws.Range("B4").Value = "5*"
$ws.Range("E4").Value = 1.6419999999999999
$ws.Range("F4").Value = 478
$ws.Range("G4").Value = 68.44
$ws.Range("H4").Value = 958

# --- Row 5: add ml_analyzed/runtime/inhibittime/numtriggers values ---
$ws.Range("E5").Value = 4.1097999999999999
$ws.Range("F5").Value = 1198
$ws.Range("G5").Value = 211.6
$ws.Range("H5").Value = 2605

# --- Row 9: fix Volume to "5*" ---
$ws.Range("B9").Value = "5*"

# --- Row 10 (new): first horz sample ---
$ws.Range("A10").Value = "D20151103T170912"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "H"
$ws.Range("D10").Value = 580
$ws.Range("E10").Value = 4.1439000000000004
$ws.Range("F10").Value = 1198
$ws.Range("G10").Value = 187
$ws.Range("H10").Value = 2403
$ws.Range("I10").Value = "First horz sample"
$ws.Range("J10").Value = "ypos moved up a bit (matlab) makes laser off more in PMTB signals"
$ws.Range("K10").Value = "9um beads, use all signals"

# --- Row 11 (new) ---
$ws.Range("A11").Value = "D20151103T173215"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "H"
$ws.Range("D11").Value = 580.1
$ws.Range("E11").Value = 4.1440999999999999
$ws.Range("F11").Value = 1198
$ws.Range("G11").Value = 188.78
$ws.Range("H11").Value = 2404
$ws.Range("I11").Value = "9um beads, use all signals"

# --- Row 12 (new) ---
$ws.Range("A12").Value = "D20151103T175435"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = "H"
$ws.Range("D12").Value = 591.6
$ws.Range("E12").Value = 4.1295999999999999
$ws.Range("F12").Value = 1198
$ws.Range("G12").Value = 190.1
$ws.Range("H12").Value = 2443
$ws.Range("I12").Value = "9um beads, use all signals"

# --- Row 13 (new) ---
$ws.Range("A13").Value = "D20151103T181654"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "H"
$ws.Range("D13").Value = 565.70000000000005
$ws.Range("E13").Value = 4.1664000000000003
$ws.Range("F13").Value = 1198
$ws.Range("G13").Value = 183.6
$ws.Range("H13").Value = 2357
$ws.Range("I13").Value = "9um beads, use all signals"

# --- Row 14 (new) ---
$ws.Range("A14").Value = "D20151103T183914"
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = "H"
$ws.Range("D14").Value = 543.4
$ws.Range("E14").Value = 4.2340999999999998
$ws.Range("F14").Value = 1198
$ws.Range("G14").Value = 181.79
$ws.Range("H14").Value = 2301
$ws.Range("I14").Value = "9um beads, use all signals"

# --- Row 15 (new) ---
$ws.Range("A15").Value = "D20151103T190133"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "H"
$ws.Range("D15").Value = 578.29999999999995
$ws.Range("E15").Value = 4.1193
$ws.Range("F15").Value = 1198
$ws.Range("G15").Value = 192.3
$ws.Range("H15").Value = 2382
$ws.Range("I15").Value = "9um beads, use all signals"

# --- Row 16 (new) ---
$ws.Range("A16").Value = "D20151103T192802"
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = "H"
$ws.Range("I16").Value = "first run switched to 4ml from 5ml"
$ws.Range("J16").Value = "9um beads, use all signals"

# move the active selection to A17, matching the post-edit cursor position
[void]$ws.Range("A17").Select()
